# Add simple mutex logic: insert a new "MUTEX_LIST" column (G) before the
# existing blank column and populate it with "m1" for the rows that
# reference a mutex.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G; existing G..K (old blank col, START, CON_NAME,
# END, INITIAL_VALUE) shift right to H..L.
$ws.Columns("G").Insert()

# Header for the new column.
$ws.Range("G1").Value = "MUTEX_LIST"

# Populate mutex usage for the relevant rows.
$ws.Range("G2").Value = "m1"
$ws.Range("G4").Value = "m1"
$ws.Range("G5").Value = "m1"

# Update the active selection as recorded in the saved workbook.
$ws.Range("G5").Select()
